$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- Row 8 ----------
$ws.Range("G8").Value = 2.95
$ws.Range("T8").Value = 9
$ws.Range("U8").Value = 15.5
$ws.Range("V8").Value = 10.5
$ws.Range("X8").Value = 26
$ws.Range("Y8").Value = 32
$ws.Range("AE8").Value = 8
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 9
$ws.Range("AI8").Value = 19.5
$ws.Range("AJ8").Value = 28

# ---------- Row 9 ----------
$ws.Range("G9").Value = 3.75
$ws.Range("AF9").Value = 8
$ws.Range("AH9").Value = 15.5
$ws.Range("AI9").Value = 17

# ---------- Row 10 ----------
$ws.Range("T10").Value = 6.4
$ws.Range("X10").Value = 24
$ws.Range("AE10").Value = 7.3
$ws.Range("AI10").Value = 28

# ---------- Row 11 ----------
$ws.Range("I11").Value = 3.75
$ws.Range("M11").Value = 2.67
$ws.Range("R11").Value = 1.87
$ws.Range("S11").Value = 1.75
$ws.Range("T11").Value = 6.4
$ws.Range("U11").Value = 8.75
$ws.Range("V11").Value = 8.5
$ws.Range("W11").Value = 17

# ---------- Row 12 (was entirely blank, now filled with odds) ----------
$ws.Range("G12").Value = 1.6
$ws.Range("H12").Value = 4.2
$ws.Range("I12").Value = 4.4
$ws.Range("N12").Value = 1.47
$ws.Range("O12").Value = 2.35
$ws.Range("R12").Value = 1.53
$ws.Range("S12").Value = 2.2
$ws.Range("T12").Value = 10
$ws.Range("U12").Value = 9.5
$ws.Range("V12").Value = 8.25
$ws.Range("W12").Value = 13
$ws.Range("X12").Value = 11.5
$ws.Range("Y12").Value = 18.5
$ws.Range("Z12").Value = 17
$ws.Range("AA12").Value = 8.75
$ws.Range("AB12").Value = 13.5
$ws.Range("AC12").Value = 45
$ws.Range("AD12").Value = 250
$ws.Range("AE12").Value = 17.5
$ws.Range("AF12").Value = 30
$ws.Range("AG12").Value = 14.5
$ws.Range("AH12").Value = 70
$ws.Range("AI12").Value = 35
$ws.Range("AJ12").Value = 32

# ---------- Row 13 ----------
$ws.Range("H13").Value = 3.3
$ws.Range("I13").Value = 2.75
$ws.Range("L13").Value = 1.36
$ws.Range("M13").Value = 3
$ws.Range("P13").Value = 1.44
$ws.Range("Q13").Value = 2.63
$ws.Range("R13").Value = 1.83
$ws.Range("S13").Value = 1.83
$ws.Range("T13").Value = 7.5
$ws.Range("Y13").Value = 34
$ws.Range("Z13").Value = 8.5
$ws.Range("AD13").Value = 301
$ws.Range("AE13").Value = 8
$ws.Range("AG13").Value = 11
$ws.Range("AH13").Value = 29

# ---------- Row 14 ----------
$ws.Range("G14").Value = 1.53
$ws.Range("H14").Value = 4.05
$ws.Range("I14").Value = 4.9
$ws.Range("J14").Value = 1.04
$ws.Range("K14").Value = 8.5
$ws.Range("L14").Value = 1.22
$ws.Range("M14").Value = 3.8
$ws.Range("N14").Value = 1.65
$ws.Range("O14").Value = 2.1
$ws.Range("P14").Value = 1.33
$ws.Range("Q14").Value = 3.05
$ws.Range("R14").Value = 1.78
$ws.Range("S14").Value = 1.93
$ws.Range("T14").Value = 7.7
$ws.Range("U14").Value = 7.7
$ws.Range("W14").Value = 11.25
$ws.Range("X14").Value = 11.75
$ws.Range("Y14").Value = 24
$ws.Range("Z14").Value = 8.5
$ws.Range("AA14").Value = 8.25
$ws.Range("AB14").Value = 16.5
$ws.Range("AC14").Value = 70
$ws.Range("AD14").Value = 500
$ws.Range("AE14").Value = 15.5
$ws.Range("AF14").Value = 30
$ws.Range("AG14").Value = 16.5
$ws.Range("AH14").Value = 90
$ws.Range("AI14").Value = 45
$ws.Range("AJ14").Value = 45
